$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest GitHub Actions scrape
$ws.Range("D2").Value = "'30.223.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "'1.861.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'235.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'0.2899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "'0.06565"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'21.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").Value = "'0.07994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "'97.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'1.860.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "'5.112"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "'0.6788"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'268.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "'30.217.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "'13.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.77%  "
$ws.Range("D19").Value = "'0.000007641"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.44%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'2.106.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'5.230"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("D24").Value = "'6.189"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'167.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'9.193"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'18.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'0.09934"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("D31").Value = "'4.339"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "'1.466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "'4.042"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'0.04717"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'1.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'0.7027"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "'2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'0.01876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'2.608"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").Value = "'6.324"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("D41").Value = "'73.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'1.939"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "'0.8400"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "'103.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'0.4144"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.171"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.055"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "'930.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'34.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "'0.05664"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
